$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "codigo electoral" (column E) values for the two data rows were numeric
# placeholders (12); replace them with the real electoral code text.
$ws.Range("E2").Value = "AST001"
$ws.Range("E3").Value = "AST001"

# Move the active selection from D4 to E4, matching the saved view state.
$ws.Range("E4").Select()

# Best-effort: nudge the saved window width to match the author's session.
$excel.ActiveWindow.Width = 11535
